$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the four new physics programs (each defaults to "Yes" in the Choose column,
# matching the existing rows' pattern).
$ws.Range("A4").Value = "RWTH_PHYSICS"
$ws.Range("B4").Value = "Yes"

$ws.Range("A5").Value = "UNI_WURZBURG_PHYSICS"
$ws.Range("B5").Value = "Yes"

$ws.Range("A6").Value = "UNI_FREIBURG_APPLIED_PHYSICS"
$ws.Range("B6").Value = "Yes"

$ws.Range("A7").Value = "KIT_PHYSICS"
$ws.Range("B7").Value = "Yes"

# Match the formatting of the existing program rows (copy A1's format down
# through the new rows, same as A2/A3 already use).
$ws.Range("A1").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("B2:B7").PasteSpecial(-4122)

# Extend the Yes/No dropdown validation (previously only B1:B3) to cover
# the new rows too.
$ws.Range("B1:B3").Validation.Delete()
$ws.Range("B1:B7").Validation.Add(3, 1, 1, """Yes,No""")
